# Drop tenure_type_label for imports and exports
#
# The "relationships" sheet (3rd sheet of the workbook) has a
# "tenure_type_label" column (D) that duplicates information already
# carried by the "tenure_type" code column (C). Remove it entirely so the
# "notes" column shifts left from E to D, and the sheet's used range
# shrinks from A1:E10 to A1:D10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Column D == "tenure_type_label" (header in D1, "Customary Rights" /
# "Water Rights" below it). Deleting the whole column shifts everything
# to its right (the "notes" column, and the blank trailing columns) one
# place to the left, and updates the sheet dimension automatically.
$ws.Columns.Item(4).Delete()
